$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 16, pushing the existing rows 16-19 down to 17-20.
# This mirrors a new week's data being added at the top of the recent-date block,
# with the old rows shifting down to make room (weekly fruit/vegetable update).
$ws.Rows.Item(16).Insert()

# Populate the newly inserted row 16 with this week's record. All of the
# descriptive columns match the surrounding rows for this market/product;
# only the date (D), volume (J) and price-per-kg (P) differ, along with the
# min/max/avg prices (K, L, M).
$ws.Cells.Item(16, 1).Value = 9
$ws.Cells.Item(16, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(16, 3).Value = "Metropolitana"
$ws.Cells.Item(16, 4).Value = 44474
$ws.Cells.Item(16, 4).Style = $ws.Cells.Item(17, 4).Style
$ws.Cells.Item(16, 4).NumberFormat = $ws.Cells.Item(17, 4).NumberFormat
$ws.Cells.Item(16, 5).Value = 13
$ws.Cells.Item(16, 6).Value = 100112010
$ws.Cells.Item(16, 7).Value = "Achicoria"
$ws.Cells.Item(16, 8).Value = "Sin especificar"
$ws.Cells.Item(16, 9).Value = "Primera"
$ws.Cells.Item(16, 10).Value = 52
$ws.Cells.Item(16, 11).Value = 5000
$ws.Cells.Item(16, 12).Value = 6000
$ws.Cells.Item(16, 13).Value = 5500
$ws.Cells.Item(16, 14).Value = "$/caja 16 unidades"
$ws.Cells.Item(16, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(16, 16).Value = 344
$ws.Cells.Item(16, 17).Value = 16
$ws.Cells.Item(16, 18).Value = "Hortaliza"
